$d = $word.ActiveDocument

# 1) Merge the split runs "2" + "8" + "/12/2022" into a single run "28/12/2022"
$d.Content.Find.Execute("28/12/2022", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "28/12/2022", 2) | Out-Null

# 2) Merge the split runs "Chess 2." + "5" into a single run "Chess 2.5"
$d.Content.Find.Execute("Chess 2.5", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Chess 2.5", 2) | Out-Null

# 3) Fill in the trailing empty list item with "Added checkmate"
$last = $d.Paragraphs.Last
$last.Range.Text = "Added checkmate"

# 4) Add a new list item after it with "Added draw"
$last = $d.Paragraphs.Last
$rng = $last.Range
$rng.Collapse(0)
$rng.InsertParagraphAfter()
$newLast = $d.Paragraphs.Last
$newLast.Range.Text = "Added draw"
